$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.954.41"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "3.514.24"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'588.62"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").Value = "'171.75"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").Value = "3.506.56"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "'0.188"
$ws.Range("E10").Value = "  -3.77%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "'0.581"
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("D13").Value = "'47.53"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").Value = "4.080.04"
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("E16").Value = "  -3.98%  "
$ws.Range("D17").Value = "'626.44"
$ws.Range("E17").Value = "  -6.07%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.523.50"
$ws.Range("E18").Value = "  -2.29%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "69.080.52"
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").Value = "'0.122"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "'17.40"
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("D22").Value = "'11.15"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").Value = "'0.887"
$ws.Range("E23").Value = "  -4.31%  "
$ws.Range("D24").Value = "'15.95"
$ws.Range("E24").Value = "  -6.33%  "
$ws.Range("D25").Value = "'97.04"
$ws.Range("E25").Value = "  -2.42%  "
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "'2.63"
$ws.Range("E28").Value = "  -4.77%  "
$ws.Range("E29").Value = "  -6.32%  "
$ws.Range("D30").Value = "'32.69"
$ws.Range("E30").Value = "  -5.07%  "
$ws.Range("D31").Value = "'8.55"
$ws.Range("E31").Value = "  -4.21%  "
$ws.Range("E32").Value = "  -5.56%  "
$ws.Range("E33").Value = "  -4.05%  "
$ws.Range("D34").Value = "'6.92"
$ws.Range("E34").Value = "  -6.05%  "
$ws.Range("D35").Value = "'639.38"
$ws.Range("E35").Value = "  +11.30%  "
$ws.Range("D36").Value = "'10.75"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("D37").Value = "'3.48"
$ws.Range("E37").Value = "  -11.23%  "
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").Value = "'57.25"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "'0.0455"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("D43").Value = "3.388.75"
$ws.Range("E43").Value = "  -5.14%  "
$ws.Range("E44").Value = "  -4.13%  "
$ws.Range("D45").Value = "'32.77"
$ws.Range("E45").Value = "  -4.83%  "
$ws.Range("D46").Value = "0.0₃0698"
$ws.Range("E46").Value = "  -4.56%  "
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("E48").Value = "  -4.63%  "
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").Value = "'132.68"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("E51").Value = "  +14.82%  "

# Re-apply default "Normal" style to cells whose numeric-looking text
# values would otherwise have been auto-converted to numbers (quote-prefix
# entry leaves a quotePrefix style flag; resetting keeps the literal text
# value while restoring the original unstyled appearance).
$ws.Range("D5,D6,D10,D12,D13,D17,D20,D21,D22,D23,D24,D25,D28,D30,D31,D34,D35,D36,D37,D39,D41,D45,D50").Style = "Normal"
